{"js": "// Eat&Reorder - Problem Statement.docx\n//\n// Three small wording tweaks (plural -> singular) in the \"Gestione acquisto\"\n// requirements list, plus the \"_GoBack\" last-edit bookmark moving from its\n// old spot (end of \"RFR4. Visualizza il catalogo degli ordini\") to the new\n// last-edited spot (right after \"...quantit\u00e0 di un\", inside the third\n// tweak), matching Word's behaviour of keeping a single _GoBack bookmark\n// that tracks the most recent edit location.\n\nconst body = context.document.body;\n\n// 1) \"Ricerca di aziende\" -> \"Ricerca di azienda\"\nlet results = body.search(\"Ricerca di aziende\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"Ricerca di azienda\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \"Ricerca di prodotti\" -> \"Ricerca di un prodotto\"\nresults = body.search(\"Ricerca di prodotti\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"Ricerca di un prodotto\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) \"Modifica della quantit\u00e0 del prodotto\" -> \"Modifica della quantit\u00e0 di un prodotto\"\nresults = body.search(\"Modifica della quantit\u00e0 del prodotto\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\n  \"Modifica della quantit\u00e0 di un prodotto\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 4) Move the \"_GoBack\" bookmark: drop the old one (after the RFR4 line) and\n//    drop a fresh, collapsed one right after \"...quantit\u00e0 di un\" (before the\n//    trailing \" prodotto\"), matching where the edit above last landed.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nresults = body.search(\"Modifica della quantit\u00e0 di un prodotto\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nconst updated = results.items[0];\n\nconst tail = updated.search(\" prodotto\", { matchCase: true });\ntail.load(\"items\");\nawait context.sync();\n\nconst insertionPoint = tail.items[0].getRange(Word.RangeLocation.start);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Eat&Reorder - Problem Statement.docx\n#\n# Three small wording tweaks (plural -> singular) in the \"Gestione acquisto\"\n# requirements list, plus the \"_GoBack\" last-edit bookmark moving from its\n# old spot (end of \"RFR4. Visualizza il catalogo degli ordini\") to the new\n# last-edited spot (right after \"...quantit\u00e0 di un\", inside the third\n# tweak) - Word keeps a single _GoBack bookmark that tracks the most recent\n# edit location, so re-adding it under its old name relocates it.\n\n$d = $word.ActiveDocument\n\n# 1) \"Ricerca di aziende\" -> \"Ricerca di azienda\"\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Ricerca di aziende\"\n$find.Execute() | Out-Null\n$range.Text = \"Ricerca di azienda\"\n\n# 2) \"Ricerca di prodotti\" -> \"Ricerca di un prodotto\"\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Ricerca di prodotti\"\n$find.Execute() | Out-Null\n$range.Text = \"Ricerca di un prodotto\"\n\n# 3) \"Modifica della quantit\u00e0 del prodotto\" -> \"Modifica della quantit\u00e0 di un prodotto\"\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Modifica della quantit\u00e0 del prodotto\"\n$find.Execute() | Out-Null\n$range.Text = \"Modifica della quantit\u00e0 di un prodotto\"\n\n# 4) Drop a fresh \"_GoBack\" bookmark right after \"...quantit\u00e0 di un\" (before\n#    the trailing \" prodotto\"), matching where the edit above last landed.\n#    Adding it under the existing name moves it off the old RFR4 spot.\n$scoped = $range.Duplicate\n$innerFind = $scoped.Find\n$innerFind.Text = \" prodotto\"\n$innerFind.Forward = $true\n$innerFind.Execute() | Out-Null\n$scoped.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $scoped) | Out-Null\n"}
